# "Reuseablity added in quick quote"
# The CreateAccount sheet's quick-quote sign-up rows (E2:E6 / F7) hold
# Selenium automation email/password test data. Refresh them with a new
# batch of generated values so the sheet can be reused for another run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("CreateAccount")

$ws.Range("E2").Value = "SeleniumVdgv@mailinator.com"
$ws.Range("E3").Value = "SeleniumsMNH@mailinator.com"
$ws.Range("E4").Value = "SeleniumfWSm@mailinator.com"
$ws.Range("E5").Value = "SeleniumTDCG@mailinator.com"
$ws.Range("E6").Value = "SeleniumSGOQ@mailinator.com"
$ws.Range("F7").Value = "Automation6512!"
